# Adds the 2022-Q3 quarter to the workbook:
#   - the current "2022-Q2" sheet is renamed to "2022-Q3" and its values are
#     overwritten with the new Q3 numbers
#   - a fresh "2022-Q2" sheet is inserted right after it, carrying the data
#     that used to live in the old "2022-Q2" sheet (unchanged)
#   - "2022-Q1" is left untouched (it just ends up 4th in tab order)
#   - the "总计" roll-up sheet gets a new top data row for 2022-Q3, pushing
#     the existing Q2/Q1 rows down one row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q2" sheet *before* touching its values, so
#    the copy keeps the old Q2 numbers. Excel places the copy right after the
#    source sheet and auto-names it "2022-Q2 (2)".
# ---------------------------------------------------------------------------
$qSrc = $wb.Worksheets.Item(2)
$qSrc.Copy($null, $qSrc)
$qCopy = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 2) Turn the original sheet into "2022-Q3" and write the new quarter's data
#    over it. Columns B, D, E, F and G hold numeric-looking text (fund codes/
#    percentages stored as strings in the source file) -- format them as text
#    first so Excel doesn't silently coerce them into numbers (which would
#    also eat leading zeros like the "003720" fund code).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3.Range("B2:B5").NumberFormat = "@"
$q3.Range("D2:G5").NumberFormat = "@"

$q3.Cells.Item(2, 2).Value = "161127"
$q3.Cells.Item(2, 3).Value = "易方达标普生物科技指数（QDII-LOF）人民币"
$q3.Cells.Item(2, 4).Value = "3.25"
$q3.Cells.Item(2, 5).Value = "94.25"
$q3.Cells.Item(2, 6).Value = "1.04"
$q3.Cells.Item(2, 7).Value = "0.0338"
$q3.Cells.Item(2, 8).Value = 4

$q3.Cells.Item(3, 2).Value = "012866"
$q3.Cells.Item(3, 3).Value = "易方达标普生物科技指数（QDII-LOF）人民币 C"
$q3.Cells.Item(3, 4).Value = "3.25"
$q3.Cells.Item(3, 5).Value = "94.25"
$q3.Cells.Item(3, 6).Value = "1.04"
$q3.Cells.Item(3, 7).Value = "0.0338"
$q3.Cells.Item(3, 8).Value = 4

$q3.Cells.Item(4, 2).Value = "003720"
$q3.Cells.Item(4, 3).Value = "易方达标普生物科技指数（QDII-LOF）美元A"
$q3.Cells.Item(4, 4).Value = "3.12"
$q3.Cells.Item(4, 5).Value = "94.25"
$q3.Cells.Item(4, 6).Value = "1.04"
$q3.Cells.Item(4, 7).Value = "0.0324"
$q3.Cells.Item(4, 8).Value = 4

# row 5's fund code (012867) is unchanged from the original sheet, only the
# remaining columns move to the new Q3 figures
$q3.Cells.Item(5, 3).Value = "易方达标普生物科技指数（QDII-LOF）美元 C"
$q3.Cells.Item(5, 4).Value = "0.13"
$q3.Cells.Item(5, 5).Value = "94.25"
$q3.Cells.Item(5, 6).Value = "1.04"
$q3.Cells.Item(5, 7).Value = "0.0014"
$q3.Cells.Item(5, 8).Value = 4

# ---------------------------------------------------------------------------
# 3) The duplicated sheet becomes the new "2022-Q2" (its data was copied
#    before we touched the original, so it still holds the old Q2 numbers).
# ---------------------------------------------------------------------------
$qCopy.Name = "2022-Q2"

# ---------------------------------------------------------------------------
# 4) Update the "总计" roll-up sheet: insert a 2022-Q3 row at the top of the
#    data and push the Q2/Q1 rows down. Row 4 is brand new, so first copy row
#    3's formatting into it (keeps the "s=2" style used by the index column)
#    before overwriting the values. Values are written bottom-up with literal
#    targets (no reading back .Value, and no Rows.Insert(), which would only
#    copy the bold header-row formatting into the new row).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Cells.Item(3, 1).Copy($total.Cells.Item(4, 1))

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q1"
$total.Cells.Item(4, 3).Value = 2
$total.Cells.Item(4, 4).Value = 0.04

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 4
$total.Cells.Item(3, 4).Value = 0.08

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.1

# ---------------------------------------------------------------------------
# 5) Keep "2022-Q1" as the selected/active tab, same as in the source file
#    (inserting the new sheets earlier in the tab order would otherwise leave
#    the old numeric tab-index pointing at the wrong sheet).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(4).Activate()
